$d = $word.ActiveDocument

# The document contains three "<id>p133v_N</id>" fragments that were each
# split across three separate runs:
#   run1: "<id>"     (Courier New, color 7f6000, sz 18)
#   run2: "p133v_N"  (color 000000, default font)
#   run3: "</id>"    (Courier New, color 7f6000, sz 18)
#
# The edit merges each triple into a single run carrying the full text
# "<id>p133v_N</id>" with the formatting of the first run (run1). Using
# Find/Replace across the contiguous text accomplishes exactly that: Word
# collapses the matched range into one run using the formatting of the
# first character of the match, and removes the now-empty extra runs.

$ids = @("p133v_1", "p133v_2", "p133v_3")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $new = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
